# Apply the "Reference" column (B) of journal citations to the RelevantBiomarkers sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contiguous row ranges sharing the same citation text (grouped by the existing
# article code already present in column A for that block of rows).
$groups = @(
    @{ Start=2;  End=6;  Text="Sachdeva et al., 2024" },
    @{ Start=7;  End=9;  Text="Zou et al., 2020" },
    @{ Start=10; End=10; Text="Amorim et al., 2022" },
    @{ Start=11; End=29; Text="Li et al., 2023" },
    @{ Start=30; End=32; Text="Yan et al., 2024" },
    @{ Start=33; End=37; Text="Lewandowicz et al., 2015" },
    @{ Start=38; End=50; Text="Nimer et al., 2023" },
    @{ Start=51; End=73; Text="Abdulwahab et al., 2019" },
    @{ Start=74; End=75; Text="Kaur et al., 2012" },
    @{ Start=76; End=78; Text="Chen et al., 2020" },
    @{ Start=79; End=82; Text="Zhao et al., 2021" },
    @{ Start=83; End=83; Text="Yu et al., 2022" },
    @{ Start=84; End=89; Text="Zhao et al., 2024" },
    @{ Start=90; End=93; Text="An et al., 2018" }
)

foreach ($g in $groups) {
    for ($r = $g.Start; $r -le $g.End; $r++) {
        $ws.Cells.Item($r, 2).Value = $g.Text
    }
}

# Widen column B to fit the new citation text and mark the whole column as
# selected, mirroring the author's manual review of the new data.
$ws.Columns.Item(2).ColumnWidth = 18.3828125
$ws.Range("B1:B1048576").Select() | Out-Null
